# Fri, Jul 31, 2020  7:07:11 PM
#
# The table on slide 16 (the cash-flow summary table, shapes:
# "Google Shape;213;p29") had its table style switched to a different
# built-in PowerPoint table style ("Medium Style 2 - Accent 1" family,
# id {05A66DC7-38E2-42DE-BDCB-3B04FCA97D36}), replacing the previous
# style id {B7E210E1-1F3B-48B9-99CE-33F94B32301E}.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)

if ($tableShape.HasTable) {
    $table = $tableShape.Table
    $table.ApplyStyle("{05A66DC7-38E2-42DE-BDCB-3B04FCA97D36}")
}
